$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (CHECKOUT_01) - swap C/D content and update E
$ws.Range("C2").Value = "User=null"
$ws.Range("D2").Value = "User=null -> gọi doGet()"
$ws.Range("E2").Value = "Lưu RedirectUrl & Chuyển Login"

# Insert a new row before the old row 3 (CHECKOUT_03), shifting CHECKOUT_03/02/04 down
$ws.Rows("3").Insert()

# New row 3 = CHECKOUT_05
$ws.Range("A3").Value = "CHECKOUT_05"
$ws.Range("B3").Value = "Vào trang Checkout (POST)"
$ws.Range("C3").Value = "User, Cart(1)"
$ws.Range("D3").Value = "Gọi doPost -> Delegated to doGet"
$ws.Range("E3").Value = "Forward Checkout.jsp"
$ws.Range("F3").Value = "OK"
$ws.Range("G3").Value = "PASS"
$ws.Range("G3").Font.Bold = $true
$ws.Range("G3").Font.Color = 32768

# Row 4 (was CHECKOUT_03, now shifted down) - swap C/D content
$ws.Range("C4").Value = "Size=0"
$ws.Range("D4").Value = "Cart size=0"

# Row 5 (was CHECKOUT_02, now shifted down) - swap C/D content
$ws.Range("C5").Value = "Cart=null"
$ws.Range("D5").Value = "User ok, Cart=null"

# Row 6 (was CHECKOUT_04, now shifted down) - update B, swap C/D content
$ws.Range("B6").Value = "Vào trang Checkout (GET)"
$ws.Range("C6").Value = "User, Cart(1)"
$ws.Range("D6").Value = "User ok, Cart ok"

# Autofit columns to match bestFit recalculated widths
$ws.Columns("A:G").AutoFit()
